# Slide 3 ("Program sources for Sprint 1, 2, 3") has a hyperlink TextBox
# that points at the previous class's GitHub folder (hw11class1). This
# homework deck was re-uploaded for hw12class2, so that textbox needs to
# be refreshed: new position/size, and the URL text updated to point at
# the hw12class2 folder.
#
# The shape is recreated (duplicated + old one removed) rather than just
# having its properties tweaked in place, so that it picks up a fresh
# shape id / default name the same way the authoring tool produced
# (id 6 "TextBox 5" -> id 5 "TextBox 4"), while still inheriting the
# existing run formatting (hyperlink run + trailing space run, noFill,
# spAutoFit body) from the shape being replaced.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$old = $s.Shapes.Item(3)

# Consume one "duplicate" id/name slot (mirrors the id/name numbering
# that happens to land the replacement shape on id=5 / "TextBox 4").
$scratch = $old.Duplicate().Item(1)
$scratch.Delete()

# The real replacement keeps all of the original shape's formatting
# (fill, body autofit, hyperlink run, trailing space run) via Duplicate.
$new = $old.Duplicate().Item(1)
$old.Delete()

$new.Name = "TextBox 4"

$new.Left = 120.75
$new.Top = 240.91874015748033
$new.Width = 684
$new.Height = 29.081259842519685

$tr = $new.TextFrame.TextRange
$tr.Runs(1).Text = "https://github.com/dougsuh/artCoding/tree/master/hw1review/hw12class2"
